# Add drop down list to data engine file
# - Adds a "Settings" sheet with lookup data for Page Name / Action Keywords / page objects
# - Inserts a "Page Name" column into "Test Steps" and fills it in for relevant rows
# - Adds defined names used by data validation dropdowns
# - Adds data validation (drop down lists) to the "Test Steps" sheet

$wb = $excel.ActiveWorkbook
$ts = $wb.Worksheets.Item("Test Steps")

# ---------------------------------------------------------------------------
# 1. Insert a new column D ("Page Name") into Test Steps; existing D/E/F shift
#    right to E/F/G.
# ---------------------------------------------------------------------------
$ts.Columns.Item(4).Insert()

$ts.Cells.Item(1, 4).Value = "Page Name"

# Login_01 block (rows 2-11)
$ts.Cells.Item(4, 4).Value = "HomePage"
$ts.Cells.Item(5, 4).Value = "LogInPage"
$ts.Cells.Item(6, 4).Value = "LogInPage"
$ts.Cells.Item(7, 4).Value = "LogInPage"
$ts.Cells.Item(9, 4).Value = "HomePage"
$ts.Cells.Item(10, 4).Value = "HomePage"

# Login_02 block (rows 12-21)
$ts.Cells.Item(14, 4).Value = "HomePage"
$ts.Cells.Item(15, 4).Value = "LogInPage"
$ts.Cells.Item(16, 4).Value = "LogInPage"
$ts.Cells.Item(17, 4).Value = "LogInPage"
$ts.Cells.Item(19, 4).Value = "HomePage"
$ts.Cells.Item(20, 4).Value = "HomePage"

# New column width (best effort; engine rounds to pixel widths)
$ts.Columns.Item(4).ColumnWidth = 18.71

# ---------------------------------------------------------------------------
# 2. The column insert does not shift the worksheet's existing hyperlinks, so
#    re-create them pointing at the now-correct column G cells.
# ---------------------------------------------------------------------------
$ts.Hyperlinks.Delete()
$ts.Hyperlinks.Add($ts.Range("G5"), "mailto:deman0590@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "deman0590@gmail.com")
$ts.Hyperlinks.Add($ts.Range("G6"), "mailto:8@aKw!7ldCyt", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "8@aKw!7ldCyt")
$ts.Hyperlinks.Add($ts.Range("G16"), "mailto:8@aKw!7ldCyt", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "8@aKw!7ldCyt")
$ts.Hyperlinks.Add($ts.Range("G15"), "mailto:deman0590@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "deman0590@gmail.com")

# ---------------------------------------------------------------------------
# 3. Add the new "Settings" sheet (placed after "Test Steps") with the
#    lookup tables used to drive the drop downs.
# ---------------------------------------------------------------------------
$wsCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($wsCount)
$settings = $wb.Worksheets.Add($null, $lastSheet)
$settings.Name = "Settings"

$settings.Cells.Item(1, 1).Value = "Action Keywords"
$settings.Cells.Item(1, 2).Value = "Page Name"
$settings.Cells.Item(1, 3).Value = "HomePage Object"
$settings.Cells.Item(1, 4).Value = "LogInPage Object"
$settings.Range("A1:D1").Font.Bold = $true

$settings.Cells.Item(2, 1).Value = "openBrowser"
$settings.Cells.Item(2, 2).Value = "HomePage"
$settings.Cells.Item(2, 3).Value = "link.login"
$settings.Cells.Item(2, 4).Value = "input.username"

$settings.Cells.Item(3, 1).Value = "navigate"
$settings.Cells.Item(3, 2).Value = "LogInPage"
$settings.Cells.Item(3, 3).Value = "button.menu"
$settings.Cells.Item(3, 4).Value = "input.password"

$settings.Cells.Item(4, 1).Value = "click"
$settings.Cells.Item(4, 3).Value = "button.logout"
$settings.Cells.Item(4, 4).Value = "button.login"

$settings.Cells.Item(5, 1).Value = "input"
$settings.Cells.Item(6, 1).Value = "waitFor"
$settings.Cells.Item(7, 1).Value = "closeBrowser"

$settings.Columns.Item(1).ColumnWidth = 16.86
$settings.Columns.Item(2).ColumnWidth = 18.14
$settings.Columns.Item(3).ColumnWidth = 18.86
$settings.Columns.Item(4).ColumnWidth = 18.43

$settings.Range("C1").Select()

# ---------------------------------------------------------------------------
# 4. Defined names that back the drop-down (Validation) lists.
# ---------------------------------------------------------------------------
$wb.Names.Add("ActionKeywords", "=Settings!`$A`$2:`$A`$30")
$wb.Names.Add("PageName", "=Settings!`$B`$2:`$B`$30")
$wb.Names.Add("HomePage", "=Settings!`$C`$2:`$C`$30")
$wb.Names.Add("LogInPage", "=Settings!`$D`$2:`$D`$30")

# ---------------------------------------------------------------------------
# 5. Data validation drop downs on the Test Steps sheet.
# ---------------------------------------------------------------------------
$ts.Range("D2:D21").Validation.Add(3, 1, 1, "=PageName")
$ts.Range("F2:F21").Validation.Add(3, 1, 1, "=ActionKeywords")
$ts.Range("E2:E21").Validation.Add(3, 1, 1, "=INDIRECT(D2)")

foreach ($dv in $ts.Validation) {
    $dv.IgnoreBlank = $true
    $dv.InCellDropdown = $true
    $dv.ShowInput = $true
    $dv.ShowError = $true
}

# ---------------------------------------------------------------------------
# 6. Restore Test Steps as the active sheet / selection.
# ---------------------------------------------------------------------------
$ts.Activate()
$ts.Range("C29").Select()
